$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 300
$ws.Range("I2").Value = 300
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 300
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -187
$ws.Range("N2").ClearContents()
$ws.Range("H29").Value = 6083.25
$ws.Range("I29").Value = 4407.3335
$ws.Range("J29").Value = 11111
$ws.Range("K29").Value = 13222.0005
$ws.Range("L29").Value = 33333
$ws.Range("M29").Value = -12941.0005
$ws.Range("N29").Value = -33895
$ws.Range("I32").Value = 886.2
$ws.Range("J32").Value = 1088.5
$ws.Range("K32").Value = 886.2
$ws.Range("L32").Value = 1088.5
$ws.Range("M32").Value = -560.2
$ws.Range("N32").Value = -1740.5
$ws.Range("H38").Value = 2659.3333
$ws.Range("J38").Value = 5750
$ws.Range("L38").Value = 17250
$ws.Range("N38").Value = -17994
$ws.Range("H39").Value = 949.7692
$ws.Range("I39").Value = 1081.5454
$ws.Range("J39").Value = 225
$ws.Range("K39").Value = 3244.6362
$ws.Range("L39").Value = 675
$ws.Range("M39").Value = -2948.6362
$ws.Range("N39").Value = -1267
$ws.Range("H40").Value = 2220.2
$ws.Range("I40").Value = 2350
$ws.Range("J40").Value = 2133.6667
$ws.Range("K40").Value = 2350
$ws.Range("L40").Value = 2133.6667
$ws.Range("M40").Value = -2175
$ws.Range("N40").Value = -2483.6667
$ws.Range("H43").Value = 550
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 550
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 550
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -688
$ws.Range("H58").Value = 1212.5454
$ws.Range("I58").Value = 42.5
$ws.Range("J58").Value = 1881.1428
$ws.Range("K58").Value = 127.5
$ws.Range("L58").Value = 5643.428400000001
$ws.Range("M58").Value = 22.5
$ws.Range("N58").Value = -5943.428400000001
$ws.Range("H64").Value = 5145.154
$ws.Range("I64").Value = 3525.125
$ws.Range("J64").Value = 7737.2
$ws.Range("K64").Value = 3525.125
$ws.Range("L64").Value = 7737.2
$ws.Range("M64").Value = -3277.125
$ws.Range("N64").Value = -8233.200000000001
$ws.Range("H67").Value = 5145.154
$ws.Range("I67").Value = 3525.125
$ws.Range("J67").Value = 7737.2
$ws.Range("K67").Value = 3525.125
$ws.Range("L67").Value = 7737.2
$ws.Range("M67").Value = -2667.125
$ws.Range("N67").Value = -9453.200000000001
$ws.Range("H106").Value = 14015167
$ws.Range("I106").Value = 16017047
$ws.Range("K106").Value = 16017047
$ws.Range("M106").Value = -16016416
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3217.2173
$ws.Range("I32").Value = 2358.795
$ws.Range("K32").Value = 2358.795
$ws.Range("M32").Value = -2071.795
$ws.Range("H61").Value = 1934.75
$ws.Range("I61").Value = 1323.625
$ws.Range("J61").Value = 4379.25
$ws.Range("K61").Value = 1323.625
$ws.Range("L61").Value = 4379.25
$ws.Range("M61").Value = -1111.625
$ws.Range("N61").Value = -4803.25
$ws.Range("H63").Value = 3491.818
$ws.Range("I63").Value = 3470.75
$ws.Range("J63").Value = 3702.5
$ws.Range("K63").Value = 3470.75
$ws.Range("L63").Value = 3702.5
$ws.Range("M63").Value = -2784.75
$ws.Range("N63").Value = -5074.5
$ws.Range("H66").Value = 3491.818
$ws.Range("I66").Value = 3470.75
$ws.Range("J66").Value = 3702.5
$ws.Range("K66").Value = 17353.75
$ws.Range("L66").Value = 18512.5
$ws.Range("M66").Value = -13921.75
$ws.Range("N66").Value = -25376.5
$ws.Range("I88").Value = 2750
$ws.Range("J88").Value = 7075
$ws.Range("K88").Value = 2750
$ws.Range("L88").Value = 7075
$ws.Range("M88").Value = -2344
$ws.Range("N88").Value = -7887
$ws.Range("I91").Value = 2750
$ws.Range("J91").Value = 7075
$ws.Range("K91").Value = 2750
$ws.Range("L91").Value = 7075
$ws.Range("M91").Value = -1346
$ws.Range("N91").Value = -9883
$ws.Range("H133").Value = 52088.555
$ws.Range("J133").Value = 52088.555
$ws.Range("L133").Value = 52088.555
$ws.Range("N133").Value = -57148.555
$ws.Range("H136").Value = 1934.75
$ws.Range("I136").Value = 1323.625
$ws.Range("J136").Value = 4379.25
$ws.Range("K136").Value = 3970.875
$ws.Range("L136").Value = 13137.75
$ws.Range("M136").Value = -1420.875
$ws.Range("N136").Value = -18237.75
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 19889.5
$ws.Range("I82").Value = 4692.75
$ws.Range("J82").Value = 50283
$ws.Range("K82").Value = 4692.75
$ws.Range("L82").Value = 50283
$ws.Range("M82").Value = -4309.75
$ws.Range("N82").Value = -51049
$ws.Range("H85").Value = 19889.5
$ws.Range("I85").Value = 4692.75
$ws.Range("J85").Value = 50283
$ws.Range("K85").Value = 4692.75
$ws.Range("L85").Value = 50283
$ws.Range("M85").Value = -3366.75
$ws.Range("N85").Value = -52935
$ws.Range("H86").Value = 16986.715
$ws.Range("I86").Value = 2200
$ws.Range("J86").Value = 22901.4
$ws.Range("K86").Value = 2200
$ws.Range("L86").Value = 22901.4
$ws.Range("M86").Value = -1077
$ws.Range("N86").Value = -25147.4
$ws.Range("H89").Value = 16986.715
$ws.Range("I89").Value = 2200
$ws.Range("J89").Value = 22901.4
$ws.Range("K89").Value = 11000
$ws.Range("L89").Value = 114507
$ws.Range("M89").Value = -5384
$ws.Range("N89").Value = -125739
$ws.Range("H94").Value = 1050.2354
$ws.Range("I94").Value = 1063.2142
$ws.Range("J94").Value = 989.6667
$ws.Range("K94").Value = 1063.2142
$ws.Range("L94").Value = 989.6667
$ws.Range("M94").Value = -612.2141999999999
$ws.Range("N94").Value = -1891.6667
$ws.Range("H105").Value = 14495537
$ws.Range("I105").Value = 18521270
$ws.Range("K105").Value = 18521270
$ws.Range("M105").Value = -18519523
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3195.5908
$ws.Range("I134").Value = 1448.9286
$ws.Range("J134").Value = 6252.25
$ws.Range("K134").Value = 4346.7858
$ws.Range("L134").Value = 18756.75
$ws.Range("M134").Value = -1811.7858
$ws.Range("N134").Value = -23826.75
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1006.23
$ws.Range("I68").Value = 773.39343
$ws.Range("J68").Value = 1370.4103
$ws.Range("K68").Value = 2320.18029
$ws.Range("L68").Value = 4111.2309
$ws.Range("M68").Value = -1509.18029
$ws.Range("N68").Value = -5733.2309
$ws.Range("H71").Value = 1006.23
$ws.Range("I71").Value = 773.39343
$ws.Range("J71").Value = 1370.4103
$ws.Range("K71").Value = 6960.54087
$ws.Range("L71").Value = 12333.6927
$ws.Range("M71").Value = -2904.54087
$ws.Range("N71").Value = -20445.6927
$ws.Range("H80").Value = 891
$ws.Range("I80").Value = 550
$ws.Range("J80").Value = 966.7778
$ws.Range("K80").Value = 1650
$ws.Range("L80").Value = 2900.3334
$ws.Range("M80").Value = -714
$ws.Range("N80").Value = -4772.3334
$ws.Range("H83").Value = 891
$ws.Range("I83").Value = 550
$ws.Range("J83").Value = 966.7778
$ws.Range("K83").Value = 4950
$ws.Range("L83").Value = 8701.0002
$ws.Range("M83").Value = -270
$ws.Range("N83").Value = -18061.0002
$ws.Range("H113").Value = 16129418
$ws.Range("I113").Value = 463.33334
$ws.Range("J113").Value = 17857520
$ws.Range("K113").Value = 1390.00002
$ws.Range("L113").Value = 53572560
$ws.Range("M113").Value = 779.9999800000001
$ws.Range("N113").Value = -53576900
$ws.Range("H129").Value = 1276.7368
$ws.Range("J129").Value = 1353.2
$ws.Range("L129").Value = 4059.6
$ws.Range("N129").Value = -14059.6
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 15000
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 15000
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 15000
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -18244
$ws.Range("H126").Value = 2065.4827
$ws.Range("I126").Value = 1692.2307
$ws.Range("J126").Value = 2368.75
$ws.Range("K126").Value = 5076.6921
$ws.Range("L126").Value = 7106.25
$ws.Range("M126").Value = -2606.6921
$ws.Range("N126").Value = -12046.25
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3294.6875
$ws.Range("I7").Value = 910
$ws.Range("J7").Value = 3453.6667
$ws.Range("K7").Value = 910
$ws.Range("L7").Value = 3453.6667
$ws.Range("M7").Value = -798
$ws.Range("N7").Value = -3677.6667
$ws.Range("H40").Value = 4788.8667
$ws.Range("I40").Value = 3500
$ws.Range("K40").Value = 3500
$ws.Range("M40").Value = -3364
$ws.Range("H46").Value = 833.5
$ws.Range("I46").Value = 900.3333
$ws.Range("J46").Value = 766.6667
$ws.Range("K46").Value = 900.3333
$ws.Range("L46").Value = 766.6667
$ws.Range("M46").Value = -712.3333
$ws.Range("N46").Value = -1142.6667
$ws.Range("H122").Value = 3990
$ws.Range("J122").Value = 3990
$ws.Range("L122").Value = 11970
$ws.Range("N122").Value = -16870
$ws.Range("H126").Value = 3294.6875
$ws.Range("I126").Value = 910
$ws.Range("J126").Value = 3453.6667
$ws.Range("K126").Value = 2730
$ws.Range("L126").Value = 10361.0001
$ws.Range("M126").Value = -260
$ws.Range("N126").Value = -15301.0001
$ws.Range("H132").Value = 5202.8667
$ws.Range("I132").Value = 3866.3333
$ws.Range("J132").Value = 6093.8887
$ws.Range("K132").Value = 11598.9999
$ws.Range("L132").Value = 18281.6661
$ws.Range("M132").Value = -9068.999899999999
$ws.Range("N132").Value = -23341.6661
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 166668340
$ws.Range("I96").Value = 200001600
$ws.Range("J96").Value = 2000
$ws.Range("K96").Value = 200001600
$ws.Range("L96").Value = 2000
$ws.Range("M96").Value = -200000227
$ws.Range("N96").Value = -4746
$ws.Range("H122").Value = 84567.336
$ws.Range("I122").Value = 143858.28
$ws.Range("K122").Value = 431574.84
$ws.Range("M122").Value = -429124.84
$ws.Range("H127").Value = 35188
$ws.Range("J127").Value = 34990
$ws.Range("L127").Value = 34990
$ws.Range("N127").Value = -44910
$ws.Range("H133").Value = 80715
$ws.Range("J133").Value = 80715
$ws.Range("L133").Value = 80715
$ws.Range("N133").Value = -90835
